$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared strings / cell text values
$ws.Range("B1").Value = "palavra"
$ws.Range("C1").Value = "correspondentes"
$ws.Range("B2").Value = "lucrou"
$ws.Range("B3").Value = "Game"

# New numeric cells in column C (rows 2-3), no special style
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0

# C1 should carry the same style as B1 (bold, bordered, centered/top aligned)
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove old row 4 entirely (A4/B4 previously held "fla"/2, now table is only 3 rows)
$ws.Rows("4:4").Delete()
